$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 1252.2142
$ws.Range("I4").Value = 1195.0769
$ws.Range("K4").Value = 1195.0769
$ws.Range("M4").Value = -1081.0769
$ws.Range("H64").Value = 4900
$ws.Range("J64").Value = 4900
$ws.Range("L64").Value = 4900
$ws.Range("N64").Value = -5396
$ws.Range("H67").Value = 4900
$ws.Range("J67").Value = 4900
$ws.Range("L67").Value = 4900
$ws.Range("N67").Value = -6616
$ws.Range("H74").Value = 146005.27
$ws.Range("I74").Value = 146005.27
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 146005.27
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = -145069.27
$ws.Range("N74").ClearContents()
$ws.Range("H77").Value = 146005.27
$ws.Range("I77").Value = 146005.27
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 730026.35
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -725346.35
$ws.Range("N77").ClearContents()
$ws.Range("H138").Value = 5837.2915
$ws.Range("J138").Value = 6788.2354
$ws.Range("L138").Value = 20364.7062
$ws.Range("N138").Value = -30644.7062

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 675
$ws.Range("I2").Value = 337.5
$ws.Range("K2").Value = 337.5
$ws.Range("M2").Value = -224.5
$ws.Range("H21").Value = 1650
$ws.Range("J21").Value = 0
$ws.Range("L21").Value = 0
$ws.Range("N21").ClearContents()
$ws.Range("H29").Value = 1000
$ws.Range("I29").Value = 1000
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 1000
$ws.Range("L29").Value = 0
$ws.Range("M29").Value = -692
$ws.Range("N29").ClearContents()
$ws.Range("H32").Value = 1951.8182
$ws.Range("I32").Value = 1951.8182
$ws.Range("K32").Value = 1951.8182
$ws.Range("M32").Value = -1664.8182
$ws.Range("H46").Value = 3997
$ws.Range("J46").Value = 3999
$ws.Range("L46").Value = 3999
$ws.Range("N46").Value = -4637
$ws.Range("H61").Value = 3668.6667
$ws.Range("I61").Value = 3668.6667
$ws.Range("K61").Value = 3668.6667
$ws.Range("M61").Value = -3456.6667
$ws.Range("H63").Value = 11465.823
$ws.Range("I63").Value = 12592.2
$ws.Range("J63").Value = 9856.714
$ws.Range("K63").Value = 12592.2
$ws.Range("L63").Value = 9856.714
$ws.Range("M63").Value = -11906.2
$ws.Range("N63").Value = -11228.714
$ws.Range("H66").Value = 11465.823
$ws.Range("I66").Value = 12592.2
$ws.Range("J66").Value = 9856.714
$ws.Range("K66").Value = 62961
$ws.Range("L66").Value = 49283.57
$ws.Range("M66").Value = -59529
$ws.Range("N66").Value = -56147.57
$ws.Range("H74").Value = 1035.7142
$ws.Range("I74").Value = 1041.6666
$ws.Range("K74").Value = 1041.6666
$ws.Range("M74").Value = -167.6666
$ws.Range("H77").Value = 1035.7142
$ws.Range("I77").Value = 1041.6666
$ws.Range("K77").Value = 5208.333000000001
$ws.Range("M77").Value = -840.3330000000005
$ws.Range("H116").Value = 675
$ws.Range("I116").Value = 337.5
$ws.Range("K116").Value = 337.5
$ws.Range("M116").Value = 1956.5
$ws.Range("H132").Value = 4434.5713
$ws.Range("I132").Value = 4434.5713
$ws.Range("K132").Value = 13303.7139
$ws.Range("M132").Value = -10773.7139
$ws.Range("H136").Value = 3668.6667
$ws.Range("I136").Value = 3668.6667
$ws.Range("K136").Value = 11006.0001
$ws.Range("M136").Value = -8456.000100000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 675
$ws.Range("I3").Value = 337.5
$ws.Range("K3").Value = 337.5
$ws.Range("M3").Value = -223.5
$ws.Range("H132").Value = 74999
$ws.Range("J132").Value = 74999
$ws.Range("L132").Value = 74999
$ws.Range("N132").Value = -85119
$ws.Range("H134").Value = 3863
$ws.Range("I134").Value = 3863
$ws.Range("K134").Value = 11589
$ws.Range("M134").Value = -9054

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2154.5
$ws.Range("I58").Value = 2060.5557
$ws.Range("K58").Value = 2060.5557
$ws.Range("M58").Value = -1857.5557
$ws.Range("H107").Value = 1943
$ws.Range("I107").Value = 2114.75
$ws.Range("K107").Value = 2114.75
$ws.Range("M107").Value = -194.75
$ws.Range("H132").Value = 3020.2222
$ws.Range("I132").Value = 3222.8572
$ws.Range("K132").Value = 9668.571599999999
$ws.Range("M132").Value = -7138.571599999999
$ws.Range("H136").Value = 2154.5
$ws.Range("I136").Value = 2060.5557
$ws.Range("K136").Value = 6181.6671
$ws.Range("M136").Value = -3631.6671

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 69965.336
$ws.Range("J37").Value = 69965.336
$ws.Range("L37").Value = 209896.008
$ws.Range("N37").Value = -210120.008
$ws.Range("H38").Value = 18.666666
$ws.Range("I38").Value = 21.285715
$ws.Range("K38").Value = 63.857145
$ws.Range("M38").Value = 283.142855

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 5141.4287
$ws.Range("I102").Value = 4998
$ws.Range("K102").Value = 4998
$ws.Range("M102").Value = -3376
$ws.Range("H113").Value = 3195
$ws.Range("I113").Value = 926.6667
$ws.Range("K113").Value = 926.6667
$ws.Range("M113").Value = 1243.3333
$ws.Range("H122").Value = 4050
$ws.Range("I122").Value = 4371.7144
$ws.Range("J122").Value = 3599.6
$ws.Range("K122").Value = 13115.1432
$ws.Range("L122").Value = 10798.8
$ws.Range("M122").Value = -10665.1432
$ws.Range("N122").Value = -15698.8
$ws.Range("H132").Value = 0
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("M132").ClearContents()
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 7997.647
$ws.Range("J7").Value = 8784.286
$ws.Range("L7").Value = 8784.286
$ws.Range("N7").Value = -9008.286
$ws.Range("H46").Value = 86550.164
$ws.Range("J46").Value = 5000
$ws.Range("L46").Value = 5000
$ws.Range("N46").Value = -5376
$ws.Range("H68").Value = 3000.5
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 3000.5
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 3000.5
$ws.Range("M68").ClearContents()
$ws.Range("N68").Value = -4498.5
$ws.Range("H71").Value = 3000.5
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 3000.5
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 15002.5
$ws.Range("M71").ClearContents()
$ws.Range("N71").Value = -22490.5
$ws.Range("H82").Value = 2390
$ws.Range("I82").Value = 2307.5
$ws.Range("K82").Value = 2307.5
$ws.Range("M82").Value = -1946.5
$ws.Range("H85").Value = 2390
$ws.Range("I85").Value = 2307.5
$ws.Range("K85").Value = 2307.5
$ws.Range("M85").Value = -1059.5
$ws.Range("H100").Value = 1934.28
$ws.Range("I100").Value = 1798.409
$ws.Range("J100").Value = 2930.6667
$ws.Range("K100").Value = 1798.409
$ws.Range("L100").Value = 2930.6667
$ws.Range("M100").Value = -1257.409
$ws.Range("N100").Value = -4012.6667
$ws.Range("H126").Value = 7997.647
$ws.Range("J126").Value = 8784.286
$ws.Range("L126").Value = 26352.858
$ws.Range("N126").Value = -31292.858
$ws.Range("H132").Value = 3596.4814
$ws.Range("I132").Value = 3786.6365
$ws.Range("K132").Value = 11359.9095
$ws.Range("M132").Value = -8829.9095

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2446.75
$ws.Range("I122").Value = 2446.75
$ws.Range("K122").Value = 7340.25
$ws.Range("M122").Value = -4890.25
$ws.Range("H135").Value = 232042.4
$ws.Range("J135").Value = 232042.4
$ws.Range("L135").Value = 232042.4
$ws.Range("N135").Value = -242182.4
$ws.Range("H138").Value = 94390
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()
